$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date-format of A10 back to the standard "date + time" format
# (it currently uses the date-only format that the new row below is about to get).
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new row of data (row 11) reported by the automatic update.
$ws.Range("A11").Value = 45856
$ws.Range("A11").NumberFormat = "YYYY-MM-DD"

$ws.Range("B11").Value = "GHT5678"
$ws.Range("C11").Value = "TESTE"
$ws.Range("D11").Value = "2025-07-18 13:58:39"
$ws.Range("E11").Value = "2025-07-18 13:58:40"
$ws.Range("F11").Value = "2025-07-18 13:58:41"
$ws.Range("G11").Value = "2025-07-18 13:58:42"
$ws.Range("H11").Value = "2025-07-18 13:58:43"
$ws.Range("I11").Value = "2025-07-18 13:58:43"
$ws.Range("J11").Value = "2025-07-18 13:58:44"
$ws.Range("K11").Value = "0:00:01"
$ws.Range("L11").Value = "0:00:01"
$ws.Range("M11").Value = "0:00:05"
$ws.Range("N11").Value = ""
$ws.Range("O11").Value = "2025-07-18 13:58:46"
$ws.Range("P11").Value = "2025-07-18 13:58:47"
$ws.Range("Q11").Value = "2025-07-18 13:58:48"
$ws.Range("R11").Value = "2025-07-18 13:58:49"
$ws.Range("S11").Value = "0:00:01"
$ws.Range("T11").Value = "0:00:01"
$ws.Range("U11").Value = "0:00:04"
$ws.Range("V11").Value = "0:00:01"
$ws.Range("W11").Value = "2025-07-18 13:58:45"
